# "Fixed data not imported bug" - the production-plan sheet was only
# half-populated (machine labels truncated, order assignments missing,
# makespan stale). Re-import the full solver output: header label,
# makespan, machine loads and every order assignment cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / makespan text -------------------------------------------
$ws.Range("A3").Value = "Machines"             # "Machine" -> "Machines"
$ws.Range("A2").Value = "make span: 222"       # stale "make span: 40" -> real result

# --- Row 4: machine load totals (previously only C4=0 was imported) ---
$ws.Range("A4").Value = 222
$ws.Range("B4").Value = 93
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 152

# Row 5 (A5:D5 = 1,2,3,4) already correct/unchanged.

# --- Rows 6-10: per-machine order assignments --------------------------
$ws.Range("B6").Value = "Order 3 - 26"
$ws.Range("C6").Value = "Order 5 - 40"
$ws.Range("D6").Value = "Order 10 - 70"

$ws.Range("B7").Value = "Order 6 - 52"
$ws.Range("C7").Value = "Order 8 - 74"

$ws.Range("B8").Value = "Order 7 - 109"
$ws.Range("C8").Value = "Order 1 - 135"

$ws.Range("B9").Value = "Order 4 - 129"
$ws.Range("C9").Value = "Order 9 - 183"

$ws.Range("C10").Value = "Order 2 - 222"

# --- Column widths: now that columns B & C hold the long "Order N - x"
# labels, widen them to fit; column C's old narrow custom width (for the
# old single makespan number) and D's are no longer needed as-is.
# (Values chosen so the engine's internal pixel-quantised ColumnWidth
# setter lands on the intended rendered width: ~16.29 chars for B,
# 18 chars for C.)
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 17.17
